# [Task 4] v2 Import excel file with street factors and improvements
#   - rename "edges sheet" -> "curves sheet"
#   - drop the leftover "Sheet1" worksheet (its shared strings go with it)
#   - update the remaining sheet's active selection

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("edges sheet")
$ws.Activate()
$ws.Range("C60").Select() | Out-Null

$ws.Name = "curves sheet"

$wb.Worksheets.Item("Sheet1").Delete() | Out-Null
